# Apply text replacements to the document per the diff.
$d = $word.ActiveDocument

$replacements = @(
    @('2025-09-23 Tuesday', '2025-09-24 Wednesday'),
    @('870×3=', '652×7='),
    @('793×2=', '776×9='),
    @('893×9=', '419×9='),
    @('567×3=', '633×6='),
    @('196×5=', '186×7='),
    @('471×9=', '377×8='),
    @('299×9=', '884×9='),
    @('479×5=', '864×3='),
    @('390×9=', '571×2='),
    @('760×7=', '273×2='),
    @('912×4=', '739×3='),
    @('576×5=', '964×4='),
    @('547×3=', '286×6='),
    @('334×3=', '960×8='),
    @('973×4=', '891×4='),
    @('205×7=', '899×4='),
    @('320×2=', '104×3='),
    @('434×2=', '146×3='),
    @('403×6=', '137×8='),
    @('894×8=', '943×6='),
    @('700×3=', '554×9='),
    @('957×7=', '281×5='),
    @('147×7=', '473×4='),
    @('760×9=', '259×4='),
    @('647×5=', '423×9='),
)

foreach ($pair in $replacements) {
    $find = $pair[0]
    $replace = $pair[1]
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2)
}

